# Add team record (Wins/Losses/Ties) columns to the NYY_2004 roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, copy formatting (bold / border / centered) from the
# existing header style used by "Unnamed: 28" (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-45) gets the team's season record repeated.
$ws.Range("AD2:AD45").Value = 101
$ws.Range("AE2:AE45").Value = 61
$ws.Range("AF2:AF45").Value = 0
